$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in A2:B5 with the new cluster counts
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 130

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 90

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 85

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 10

# Remove row 6 entirely (was A6=4, B6=35) so the used range shrinks to A1:B5
$ws.Range("A6:B6").Delete()
